$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added to the series. Insert a new
# row at position 10 (pushing the existing row 10 and everything below
# it down by one row) and populate it with the new record.
$ws.Rows.Item(10).Insert()

# Carry over the constant/static columns from the (now shifted) row 11,
# which still holds the same market/category/quality/unit metadata.
$ws.Range("A10:C10").Value2 = $ws.Range("A11:C11").Value2
$ws.Range("E10:I10").Value2 = $ws.Range("E11:I11").Value2
$ws.Range("N10:O10").Value2 = $ws.Range("N11:O11").Value2
$ws.Range("Q10:R10").Value2 = $ws.Range("Q11:R11").Value2

# New row's own data (date + volume/price figures).
$ws.Cells.Item(10, 4).Value = "2022-12-06"
$ws.Cells.Item(10, 10).Value = 300
$ws.Cells.Item(10, 11).Value = 900
$ws.Cells.Item(10, 12).Value = 1000
$ws.Cells.Item(10, 13).Value = 950
$ws.Cells.Item(10, 16).Value = 475
